# Natmi following Dr Hou advice
#
# The LR-pairs_lrc2p/Il1b-Il1r2 sheet is rebuilt with the (now 3-cluster:
# ECs, FAPs, sCs) specificity values recomputed, producing a full
# sender x receiver grid (6 data rows instead of the original 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="ECs";  B="Il1b"; C="Il1r2"; D="ECs";  E=3; F=1; G="2105.707937";          H="6317.123811";          I="0.9998149610961508";   J="0.9998149610961509";   K=3; L=1; M="373.8120643333334";  N="1121.436193";  O="0.9916504392533645"; P="0.9916504392533645"; Q="787139.0308130546";   R="7084251.277317492";   S="0.9914669453430834";    T="0.9914669453430835" },
    @{ Row=3; A="ECs";  B="Il1b"; C="Il1r2"; D="FAPs"; E=3; F=1; G="2105.707937";          H="6317.123811";          I="0.9998149610961508";   J="0.9998149610961509";   K=3; L=1; M="3.147446333333333";   N="9.442339";      O="0.00834956074663552"; P="0.00834956074663552"; Q="6627.602725381548";   R="59648.42452843393";   S="0.00834801575306734";   T="0.00834801575306734" },
    @{ Row=4; A="FAPs"; B="Il1b"; C="Il1r2"; D="ECs";  E=1; F="0.3333333333333333"; G="0.3524213333333333"; H="1.057264";    I="0.0001673338051705887"; J="0.0001673338051705887"; K=3; L=1; M="373.8120643333334";  N="1121.436193";  O="0.9916504392533645"; P="0.9916504392533645"; Q="131.7393461284391";  R="1185.654115155952";  S="0.0001659366413993512"; T="0.0001659366413993512" },
    @{ Row=5; A="FAPs"; B="Il1b"; C="Il1r2"; D="FAPs"; E=1; F="0.3333333333333333"; G="0.3524213333333333"; H="1.057264";    I="0.0001673338051705887"; J="0.0001673338051705887"; K=3; L=1; M="3.147446333333333";   N="9.442339";      O="0.00834956074663552"; P="0.00834956074663552"; Q="1.109227233388444";  R="9.983045100496";     S="0.000001397163771237503"; T="0.000001397163771237503" },
    @{ Row=6; A="sCs";  B="Il1b"; C="Il1r2"; D="ECs";  E=1; F="0.3333333333333333"; G="0.03728866666666666"; H="0.111866"; I="0.00001770509867848813"; J="0.00001770509867848813"; K=3; L=1; M="373.8120643333334";  N="1121.436193";  O="0.9916504392533645"; P="0.9916504392533645"; Q="13.93895346290422";  R="125.450581166138";   S="0.00001755726888154692"; T="0.00001755726888154692" },
    @{ Row=7; A="sCs";  B="Il1b"; C="Il1r2"; D="FAPs"; E=1; F="0.3333333333333333"; G="0.03728866666666666"; H="0.111866"; I="0.00001770509867848813"; J="0.00001770509867848813"; K=3; L=1; M="3.147446333333333";   N="9.442339";      O="0.00834956074663552"; P="0.00834956074663552"; Q="0.1173640771748889"; R="1.056276694574";     S="0.0000001478297969412129"; T="0.0000001478297969412129" }
)

$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    $rowIndex = $r["Row"]
    $colNum = 1
    foreach ($col in $colLetters) {
        $ws.Cells.Item($rowIndex, $colNum).Value2 = $r[$col]
        $colNum = $colNum + 1
    }
}
